$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (the default, unstyled, data-row style) to stamp on
# "empty but present" cells so they materialize in the sheet without
# picking up a new/different style.
$defaultStyle = $ws.Range("A2").Style

# --- Row 2: update existing task entry ---
$ws.Range("A2").Value = "sbkuzh"
$ws.Range("B2").Value = "task_2024-07-19_ZBS_SMALL_ERROR"
$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = "2024-07-19"

# --- Row 3: new ERROR task entry ---
$ws.Range("A3").Value = "sbkzhk"
$ws.Range("B3").Value = "task_2024-07-19_ZBS_SMALL_ERROR"
$ws.Range("C3").Style = $defaultStyle
$ws.Range("D3").Style = $defaultStyle
$ws.Range("E3").Style = $defaultStyle
$ws.Range("F3").NumberFormat = "@"
$ws.Range("F3").Value = "2024-07-19"
$ws.Range("G3").Value = "SMALL"
$ws.Range("H3").Value = "ERROR"
$ws.Range("I3").Value = 'task_2024-07-19_ZBS_SMALL_DONE: workflow broken, new tasks must have "NEW" state and not "DONE"'

# --- Row 4: new DONE task entry (start date of the task added) ---
$ws.Range("A4").Value = "sbkzbs"
$ws.Range("B4").Value = "task_2024-07-19_ZBS1_SMALL_DONE"
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2024-07-19 23:40:59"
$ws.Range("D4").Style = $defaultStyle
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2024-07-19 23:41:06"
$ws.Range("F4").NumberFormat = "@"
$ws.Range("F4").Value = "2024-07-19"
$ws.Range("G4").Value = "SMALL"
$ws.Range("H4").Value = "DONE"
$ws.Range("I4").Style = $defaultStyle
